$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 244, shifting the existing
# rows 244:251 down to 245:252.
$ws.Rows.Item(244).Insert()

# Populate the new row 244 with the new data record.
$ws.Cells.Item(244, 1).Value = 6
$ws.Cells.Item(244, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(244, 3).Value = "Metropolitana"
$ws.Cells.Item(244, 4).Value = 44509
$ws.Cells.Item(244, 4).NumberFormat = $ws.Cells.Item(245, 4).NumberFormat
$ws.Cells.Item(244, 5).Value = 13
$ws.Cells.Item(244, 6).Value = 100112032
$ws.Cells.Item(244, 7).Value = "Zapallo italiano"
$ws.Cells.Item(244, 8).Value = "Sin especificar"
$ws.Cells.Item(244, 9).Value = "Primera"
$ws.Cells.Item(244, 10).Value = 1400
$ws.Cells.Item(244, 11).Value = 5000
$ws.Cells.Item(244, 12).Value = 6000
$ws.Cells.Item(244, 13).Value = 5536
$ws.Cells.Item(244, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(244, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(244, 16).Value = 111
$ws.Cells.Item(244, 17).Value = 50
$ws.Cells.Item(244, 18).Value = "Hortaliza"
